# The template previously had three Jinja-style docxtpl tags
# ({%p if ... %} / {%p for ... %}) whose text had been split across
# several <w:r> runs by Word's spell-checker, with <w:proofErr>
# spellStart/spellEnd markers wrapped around the "suspicious" bare
# variable name in the middle (e.g. "pending_actions_between_parties_yes").
# That run-splitting caused a rendering/formatting glitch for the
# protectee lists. The fix is purely textual: re-typing each tag's
# text as one unbroken string merges the runs/removes the proofErr
# markers without touching anything else on the paragraph.
#
# Doing a Find/Replace of the exact same text (old phrase -> new phrase,
# where both happen to read identically) forces Word to rewrite that
# stretch of the paragraph as a single contiguous run, exactly mirroring
# the collapsed <w:r>/<w:t> the diff expects.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "{%p if pending_actions_between_parties_yes %}", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "{%p if pending_actions_between_parties_yes %}", 2)

$d.Content.Find.Execute(
    "{%p if orders_judgments_re_parties_yes %}", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "{%p if orders_judgments_re_parties_yes %}", 2)

$d.Content.Find.Execute(
    "{%p for order in orders_judgments %}", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "{%p for order in orders_judgments %}", 2)
